$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Removing less than USD 5 price from extrapolation calibration because it is just a noise
# Updated recalculated ABSM1_RN / M1_RN / CM2_RN / CMN3_RN / CMN4_RN values for the affected rows
$ws.Range("D3").Value = 120402.5964248481
$ws.Range("E3").Value = 0.01905326443683382
$ws.Range("F3").Value = 0.1983945371478865
$ws.Range("G3").Value = -1.589223172820981
$ws.Range("H3").Value = 17.86714729572637
$ws.Range("D5").Value = 121968.407769422
$ws.Range("E5").Value = -0.003556305241997975
$ws.Range("F5").Value = 0.2353462284745418
$ws.Range("G5").Value = -1.131481576085093
$ws.Range("H5").Value = 12.18325821204666
$ws.Range("D7").Value = 123282.7130233925
$ws.Range("E7").Value = -0.01743438741688789
$ws.Range("F7").Value = 0.2554441472190375
$ws.Range("G7").Value = -1.15314090507522
$ws.Range("H7").Value = 11.14723771527033
$ws.Range("D8").Value = 124755.3145441186
$ws.Range("E8").Value = -0.03588622783839419
$ws.Range("F8").Value = 0.208740667019779
$ws.Range("G8").Value = -0.7266885209072331
$ws.Range("H8").Value = 6.656960766096915
$ws.Range("D9").Value = 126176.5369015856
$ws.Range("E9").Value = -0.07246590519619568
$ws.Range("F9").Value = 0.3449171757093097
$ws.Range("G9").Value = -1.728044946501899
$ws.Range("H9").Value = 10.97324724333234
$ws.Range("D10").Value = 127769.1945475637
$ws.Range("E10").Value = -0.1079304521415401
$ws.Range("F10").Value = 0.4398238573457138
$ws.Range("G10").Value = -1.939634743921901
$ws.Range("H10").Value = 9.973633813997481
$ws.Range("D11").Value = 129908.7980415271
$ws.Range("E11").Value = -0.1859665212257841
$ws.Range("F11").Value = 0.7767202762314407
$ws.Range("G11").Value = -2.641894920643084
$ws.Range("H11").Value = 12.968103125291
$ws.Range("D13").Value = 120419.1111750395
$ws.Range("E13").Value = 0.1052671926194933
$ws.Range("F13").Value = 0.1490563546392251
$ws.Range("G13").Value = -0.5495603700888898
$ws.Range("H13").Value = 10.46839658534563
$ws.Range("D15").Value = 120448.4228209491
$ws.Range("E15").Value = 0.09858012152733832
$ws.Range("F15").Value = 0.1511755944397735
$ws.Range("G15").Value = -0.5728183654219944
$ws.Range("H15").Value = 11.61092752899345
$ws.Range("D17").Value = 120350.8024802794
$ws.Range("E17").Value = 0.08082528932054393
$ws.Range("F17").Value = 0.1517672747021148
$ws.Range("G17").Value = -0.8535885794979272
$ws.Range("H17").Value = 7.684185731442296
$ws.Range("D18").Value = 120431.5219054105
$ws.Range("E18").Value = 0.05765950261233897
$ws.Range("F18").Value = 0.1607735371443851
$ws.Range("G18").Value = -0.6256348215696956
$ws.Range("H18").Value = 9.098686655553168
